$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 812; all existing rows 812..853 shift down to 813..854
$ws.Rows.Item(812).Insert()

# Write B/C/D as normal literal values
$ws.Range("B812").Value = "木"
$ws.Range("C812").Value = 22
$ws.Range("D812").Value = 201

# A812 must hold the literal text "2026/02/12" (not get auto-converted into a date
# serial number the way a plain .Value assignment of a date-like string would be).
# Enter it as a formula that evaluates to the text string, then collapse the
# formula down to a plain value via copy / paste-values so the cell ends up as a
# plain text cell, matching the rest of the date column.
$ws.Range("A812").Formula = "=""2026/02/12"""
$ws.Range("A812").Copy()
$ws.Range("A812").PasteSpecial(-4163)
